$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row for rocksdb database entry
$ws.Range("A6").Value = 10002
$ws.Range("B6").Value = "rocksdb"

# Column C header changed from "IP" to "Address"
$ws.Range("C4").Value = "Address"

$ws.Range("C6").Value = "/data/db"

# Reflect the final selected cell as left by the author when saving
$ws.Range("I31").Select() | Out-Null
